$wb = $excel.ActiveWorkbook

# --- Update the summary text on "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Range("A1")
$oldText = $cellA1.Value2
$newText = $oldText -replace [regex]::Escape("✅ 1000 Bs = 10.36 = 41710.88 pesos"), "✅ 1000 Bs = 9.9 = 39852.27 pesos"
$newText = $newText -replace [regex]::Escape("✅ 41710.88 pesos = 10.32 = 967.4 Bs"), "✅ 39852.27 pesos = 9.82 = 959.01 Bs"
$cellA1.Value = $newText

# --- Update the rate cells on "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 100.998
$ws2.Range("O10").Value = 4025
$ws2.Range("N12").Value = 4060
$ws2.Range("O12").Value = 97.7
